# "Generate Report for handback"
# The row for file acb8a402-b60a-4eac-b7c6-6347f3b08d95.md has been handed
# back and is now in sync with en-US. Update the Overview sheet's status
# columns plus the per-locale sheets' Status and Latest Handback DateTime
# columns to reflect the handback.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Handed back: in sync with en-US"
$overview.Range("C3").Value = "Handed back: in sync with en-US"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B3").Value = "Handed back: in sync with en-US"
$zhcn.Range("G3").Value = "2016-02-15 03:47:18"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B3").Value = "Handed back: in sync with en-US"
$dede.Range("G3").Value = "2016-02-15 03:47:42"
